# Weekly fruit/vegetable data update: a new price record (week of 2021-11-23)
# is inserted as row 158, pushing the previously-existing rows 158-240 down
# to 159-241 (the sheet's used range grows from A1:R240 to A1:R241).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 158, shifting rows 158:240 down to 159:241.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new observation.
$ws.Cells.Item(158, 1).Value  = 9
$ws.Cells.Item(158, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(158, 3).Value  = "Metropolitana"
$ws.Cells.Item(158, 4).Value  = 44523
$ws.Cells.Item(158, 5).Value  = 13
$ws.Cells.Item(158, 6).Value  = 100112032
$ws.Cells.Item(158, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(158, 8).Value  = "Sin especificar"
$ws.Cells.Item(158, 9).Value  = "Primera"
$ws.Cells.Item(158, 10).Value = 79
$ws.Cells.Item(158, 11).Value = 7000
$ws.Cells.Item(158, 12).Value = 8000
$ws.Cells.Item(158, 13).Value = 7494
$ws.Cells.Item(158, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(158, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(158, 16).Value = 150
$ws.Cells.Item(158, 17).Value = 50
$ws.Cells.Item(158, 18).Value = "Hortaliza"
